$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Extend the used range down to row 23 with the same bordered style as the
#    existing data rows, so every new cell already carries style index 1 before
#    we start writing values into it.
$styleSource = $ws.Range("A11:F11")
$styleSource.Copy()
$ws.Range("A15:F23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Columns A, D and E hold numeric-looking text (IDs / qty / sequence numbers)
#    that must stay stored as text (shared strings), exactly like the original
#    workbook. Marking the columns as Text *before* assigning the values stops
#    Excel from re-interpreting them as numbers.
$ws.Range("A2:A23").NumberFormat = "@"
$ws.Range("D2:D23").NumberFormat = "@"
$ws.Range("E2:E23").NumberFormat = "@"

# 3) Write the final values row by row.
# Row 2
$ws.Range("A2").Value = "20024079"
$ws.Range("B2").Value = "TELUR AYM NEGERI BKL"
$ws.Range("C2").Value = "RPROCL"
$ws.Range("D2").Value = "1"
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "PT,(E-1H)"

# Row 3
$ws.Range("A3").Value = "20033350"
$ws.Range("B3").Value = "G/F FRS MLK F/CRM320"
$ws.Range("C3").Value = "RPROCL"
$ws.Range("D3").Value = "3"
$ws.Range("E3").Value = "23"
$ws.Range("F3").Value = "RT,(E-3H)"

# Row 4
$ws.Range("A4").Value = "20088719"
$ws.Range("B4").Value = "KNZLER SNGLES ORG 65"
$ws.Range("C4").Value = "RPROCL"
$ws.Range("D4").Value = "3"
$ws.Range("E4").Value = "36"
$ws.Range("F4").Value = "RT,(E-7H)"

# Row 5
$ws.Range("A5").Value = "20088730"
$ws.Range("B5").Value = "KNZLER SNGLES KJU 60"
$ws.Range("C5").Value = "RPROCL"
$ws.Range("D5").Value = "3"
$ws.Range("E5").Value = "37"
$ws.Range("F5").Value = "RT,(E-7H)"

# Row 6
$ws.Range("A6").Value = "20106308"
$ws.Range("B6").Value = "KNZLER SNGLES HOT 65"
$ws.Range("C6").Value = "RPROCL"
$ws.Range("D6").Value = "3"
$ws.Range("E6").Value = "39"
$ws.Range("F6").Value = "RT,(E-7H)"

# Row 7
$ws.Range("A7").Value = "20126489"
$ws.Range("B7").Value = "KNZLER SNGL GCHJNG60"
$ws.Range("C7").Value = "RPROCL"
$ws.Range("D7").Value = "3"
$ws.Range("E7").Value = "91"
$ws.Range("F7").Value = "RT,(E-7H)"

# Row 8
$ws.Range("A8").Value = "20137459"
$ws.Range("B8").Value = "KZLR BAKSO GCHJNG 55"
$ws.Range("C8").Value = "RPROCL"
$ws.Range("D8").Value = "3"
$ws.Range("E8").Value = "126"
$ws.Range("F8").Value = "RT,(E-7H)"

# Row 9
$ws.Range("A9").Value = "20137391"
$ws.Range("B9").Value = "SG SOSIS JPG HOT 55G"
$ws.Range("C9").Value = "RPROCL"
$ws.Range("D9").Value = "3"
$ws.Range("E9").Value = "128"
$ws.Range("F9").Value = "RT,(E-7H)"

# Row 10
$ws.Range("A10").Value = "20137392"
$ws.Range("B10").Value = "SG SOSIS JPG KEJU 55"
$ws.Range("C10").Value = "RPROCL"
$ws.Range("D10").Value = "3"
$ws.Range("E10").Value = "129"
$ws.Range("F10").Value = "RT,(E-7H)"

# Row 11
$ws.Range("A11").Value = "20138147"
$ws.Range("B11").Value = "CMORY EAT MLK HZL 80"
$ws.Range("C11").Value = "RPROCL"
$ws.Range("D11").Value = "3"
$ws.Range("E11").Value = "139"
$ws.Range("F11").Value = "RT,(E-7H)"

# Row 12
$ws.Range("A12").Value = "20138150"
$ws.Range("B12").Value = "CMORY EAT MLK CHO 80"
$ws.Range("C12").Value = "RPROCL"
$ws.Range("D12").Value = "3"
$ws.Range("E12").Value = "140"
$ws.Range("F12").Value = "RT,(E-7H)"

# Row 13
$ws.Range("A13").Value = "20138151"
$ws.Range("B13").Value = "CMORY EAT MLK MRIE80"
$ws.Range("C13").Value = "RPROCL"
$ws.Range("D13").Value = "3"
$ws.Range("E13").Value = "141"
$ws.Range("F13").Value = "RT,(E-7H)"

# Row 14
$ws.Range("A14").Value = "20139684"
$ws.Range("B14").Value = "FIESTA BAKSO MN.LAVA"
$ws.Range("C14").Value = "RPROCL"
$ws.Range("D14").Value = "3"
$ws.Range("E14").Value = "154"
$ws.Range("F14").Value = "RT,(E-7H)"

# Row 15
$ws.Range("A15").Value = "20138476"
$ws.Range("B15").Value = "FIESTA S.RTG H.BBQ60"
$ws.Range("C15").Value = "RPROCL"
$ws.Range("D15").Value = "3"
$ws.Range("E15").Value = "158"
$ws.Range("F15").Value = "RT,(E-7H)"

# Row 16
$ws.Range("A16").Value = "20139811"
$ws.Range("B16").Value = "MABELL SS S.TEMPONG"
$ws.Range("C16").Value = "RPROCL"
$ws.Range("D16").Value = "3"
$ws.Range("E16").Value = "161"
$ws.Range("F16").Value = "RT,(E-7H)"

# Row 17
$ws.Range("A17").Value = "20139812"
$ws.Range("B17").Value = "MABELL GOCHUJANG 60G"
$ws.Range("C17").Value = "RPROCL"
$ws.Range("D17").Value = "3"
$ws.Range("E17").Value = "162"
$ws.Range("F17").Value = "RT,(E-7H)"

# Row 18
$ws.Range("A18").Value = "10037636"
$ws.Range("B18").Value = "CMPN HULA KCG.HJU 45"
$ws.Range("C18").Value = "RPROCL"
$ws.Range("D18").Value = "4"
$ws.Range("E18").Value = "30"
$ws.Range("F18").Value = "RT"

# Row 19
$ws.Range("A19").Value = "20026370"
$ws.Range("B19").Value = "CMPN HULA TAP KTN 45"
$ws.Range("C19").Value = "RPROCL"
$ws.Range("D19").Value = "4"
$ws.Range("E19").Value = "34"
$ws.Range("F19").Value = "RT"

# Row 20
$ws.Range("A20").Value = "20113120"
$ws.Range("B20").Value = "AICE CHOCO ALMOND 90"
$ws.Range("C20").Value = "RPROCL"
$ws.Range("D20").Value = "4"
$ws.Range("E20").Value = "82"
$ws.Range("F20").Value = "RT,(E-3B)"

# Row 21
$ws.Range("A21").Value = "20131384"
$ws.Range("B21").Value = "CMPNA HULA ALPKAT 45"
$ws.Range("C21").Value = "RPROCL"
$ws.Range("D21").Value = "4"
$ws.Range("E21").Value = "173"
$ws.Range("F21").Value = "RT,(E-1B)"

# Row 22
$ws.Range("A22").Value = "20134511"
$ws.Range("B22").Value = "AICE CLSC CHO ALMD90"
$ws.Range("C22").Value = "RPROCL"
$ws.Range("D22").Value = "4"
$ws.Range("E22").Value = "179"
$ws.Range("F22").Value = "RT,(E-1B)"

# Row 23
$ws.Range("A23").Value = "20140183"
$ws.Range("B23").Value = "WALLS X PC PSL 90ML"
$ws.Range("C23").Value = "RPROCL"
$ws.Range("D23").Value = "4"
$ws.Range("E23").Value = "240"
$ws.Range("F23").Value = "RT,(E-1B)"

# 4) Restore the original (General / bordered, style index 1) number format on the
#    numeric-text columns now that the values are locked in as text -- this keeps
#    the cell style identical to the rest of the sheet.
$generalFormat = $ws.Range("B2")
$generalFormat.Copy()
$ws.Range("A2:A23").PasteSpecial(-4122)
$ws.Range("D2:D23").PasteSpecial(-4122)
$ws.Range("E2:E23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

